$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-08 Monday" "2025-12-09 Tuesday"
Replace-Text "12×61=" "51×71="
Replace-Text "48×56=" "49×84="
Replace-Text "43×58=" "38×30="
Replace-Text "92×39=" "56×55="
Replace-Text "63×61=" "56×43="
Replace-Text "23×44=" "24×76="
Replace-Text "82×14=" "72×24="
Replace-Text "23×91=" "53×68="
Replace-Text "63×95=" "55×81="
Replace-Text "27×26=" "71×43="
Replace-Text "55×83=" "21×36="
Replace-Text "93×13=" "87×70="
Replace-Text "71×59=" "38×81="
Replace-Text "37×39=" "70×57="
Replace-Text "59×12=" "71×67="
Replace-Text "38×38=" "40×69="
Replace-Text "64×80=" "23×81="
Replace-Text "35×25=" "71×77="
Replace-Text "16×23=" "91×77="
Replace-Text "12×62=" "68×14="
Replace-Text "26×36=" "49×49="
Replace-Text "45×95=" "47×87="
Replace-Text "71×73=" "91×88="
Replace-Text "61×17=" "98×73="
Replace-Text "75×21=" "88×24="
